$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.214.90"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.293.02"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'586.09"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").Value = "'180.91"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "'0.653"
$ws.Range("E7").Value = "  +9.43%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "'6.76"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").Value = "'0.406"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").Value = "3.866.75"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("D14").Value = "66.239.67"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'26.50"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.317.98"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000164"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "'437.13"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'13.29"
$ws.Range("D20").Value = "'5.51"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").Value = "'7.46"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").Value = "'72.51"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'5.69"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "3.426.07"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "'0.512"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'0.197"
$ws.Range("E27").Value = "  +3.79%  "
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").Value = "'8.87"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "'22.38"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'5.22"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "'6.64"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").Value = "'158.20"
$ws.Range("D38").Value = "'1.42"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("D39").Value = "'26.59"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("D40").Value = "'1.78"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("D41").Value = "2.800.33"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").Value = "'0.775"
$ws.Range("D43").Value = "'4.36"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "'40.30"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "'6.11"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'0.0663"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").Value = "'2.32"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "'320.60"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "'23.35"
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("D50").Value = "'0.0268"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  +6.65%  "
